# Insert a new "Force" config row into the GuildConfig sheet, between the
# existing "Ref" row (row 7) and the "Upload" row (row 8), pushing the
# "Upload" row and everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 (shifts old rows 8-12 down to 9-13,
# and automatically extends dependent ranges such as the dataValidations
# that covered A7:A8 / B7:J8 to A7:A9 / B7:J9).
$ws.Rows(8).Insert()

# Copy the formatting of the row above (the "Ref" row, which has the same
# look every boolean-flag row in this block uses) onto the new row so the
# new cells carry the correct styles (label cell vs. boolean-value cells).
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G8").PasteSpecial(-4122)

# Populate the new row's content: label + six FALSE boolean flags, matching
# the pattern of the surrounding rows (Ref/Upload).
$ws.Range("A8").Value = "Force"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $false
$ws.Range("G8").Value = $false
